# "half of the amazon data." — add a new "accuracy RATE mean" column
# (= accuracy mean / 69) between the existing "accuracy mean" (D) and
# "accuracy stdev" (old E, now shifted to F) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column at E; this pushes the old "accuracy stdev"
# column (and its data) from E to F automatically.
$ws.Columns("E:E").Insert()

# Header for the new column.
$ws.Range("E1").Value = "accuracy RATE mean"

# Row 2 was authored as a standalone formula...
$ws.Range("E2").Formula = "=`$D2/69"
# ...while rows 3-15 were filled/copied together, so they share one
# formula group (si="0", master on E3, ref="E3:E15").
$ws.Range("E3:E15").Formula = "=`$D3/69"

# Column widths that got set (manually or via autofit) after the insert.
$ws.Columns("D:D").ColumnWidth = 13.666666666666666
$ws.Columns("E:E").ColumnWidth = 15.330729166666666

# Final selection left on the sheet.
$ws.Range("D19").Select()
$excel.ActiveWindow.ScrollRow = 3
